$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.459612070389937, 1.667794583268128, 26.21740644021617, 0.496779210170732, 29.84159230404497)
    3 = @(0.127881588408715, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 0.8245905782990299)
    4 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    5 = @(0.127881588408715, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 4.837881874639075)
    6 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    7 = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
